$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.019.52"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "'2.794.32"
$ws.Range("E3").Value = "  -1.64%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'359.71"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").Value = "'109.81"
$ws.Range("E6").Value = "  -3.04%  "

$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "  -2.69%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("D10").Value = "'40.28"
$ws.Range("E10").Value = "  -3.16%  "

$ws.Range("E11").Value = "  -1.45%  "

$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").Value = "'19.53"
$ws.Range("E13").Value = "  -2.58%  "

$ws.Range("D14").Value = "'7.60"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").Value = "'3.231.54"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "'2.804.11"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "'0.944"
$ws.Range("E17").Value = "  +3.89%  "

$ws.Range("D18").Value = "'51.961.09"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("D20").Value = "'3.12"
$ws.Range("E20").Value = "  -1.17%  "

$ws.Range("D21").Value = "'13.14"
$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("D22").Value = "'0.0₃0979"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'270.58"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'70.26"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").Value = "'26.57"
$ws.Range("E26").Value = "  -2.21%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  +15.31%  "

$ws.Range("D29").Value = "'10.33"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").Value = "'2.28"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").Value = "'0.0474"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").Value = "'52.22"
$ws.Range("E32").Value = "  -3.92%  "

$ws.Range("D33").Value = "'34.52"
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").Value = "'0.0849"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("E36").Value = "  -5.62%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "'19.10"
$ws.Range("E38").Value = "  +3.50%  "

$ws.Range("D39").Value = "'3.22"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").Value = "'2.00"
$ws.Range("E40").Value = "  -3.68%  "

$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("E42").Value = "  -2.12%  "

$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").Value = "'119.71"
$ws.Range("E44").Value = "  -6.61%  "

$ws.Range("D45").Value = "'21.91"
$ws.Range("E45").Value = "  -8.81%  "

$ws.Range("D46").Value = "'2.088.91"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("E47").Value = "  -4.43%  "

$ws.Range("D49").Value = "'5.76"
$ws.Range("E49").Value = "  -2.22%  "

$ws.Range("D50").Value = "'0.958"
$ws.Range("E50").Value = "  -4.36%  "

$ws.Range("D51").Value = "'8.85"
$ws.Range("E51").Value = "  -2.61%  "
